$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @("43-15=","77+14=","85-28=","98-59=","13+28=","9+8=","94-86=","96-37=","80-51=","7+68=","26+6=","28+58=","29+22=","47+4=","50-7=","24+68=","42-38=","47+8=","42+39=","46-19=","25+29=","72-35=","75-67=","91-16=","7+14=","86-37=","70-6=","59+36=","60-29=","33+19=","27+64=","91-27=","84-78=","77-19=","65-37=","31-24=","40-12=","11-2=","44-6=","47+27=","9+64=","95-6=","18+39=","8+63=","65-18=","12+49=","38+13=","4+58=","60-47=","71-28=","80-28=","36+57=","32+39=","4+67=","7+54=","73-17=","94-75=","50-7=","55+38=","90-32=","29+56=","61-32=","5+29=","16+58=","80-13=","7+34=","49+14=","54-38=","13+68=","16-9=","17+79=","6+56=","9+59=","43-38=","49+9=","51-26=","47+44=","9+38=","84+7=","33+19=","54-7=","9+24=","49+9=","36+18=","54-25=","73-57=","68+14=","18+27=","19+19=","90-25=","12-6=","48+25=","87-49=","19+62=","31-23=","13+9=","74-28=","12+69=","89+2=","46-19=")

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"
